# Updated symbol list - apply latest price/volume values to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref="D2"; Value="291.50"}
    @{Ref="E2"; Value="-3.75%"}
    @{Ref="E3"; Value="-3.80%"}
    @{Ref="D4"; Value="4.868"}
    @{Ref="E4"; Value="-2.54%"}
    @{Ref="D5"; Value="0.07209"}
    @{Ref="E5"; Value="-8.88%"}
    @{Ref="D6"; Value="7.668"}
    @{Ref="E6"; Value="-2.47%"}
    @{Ref="D7"; Value="1.728"}
    @{Ref="E7"; Value="-17.62%"}
    @{Ref="D8"; Value="3.763"}
    @{Ref="E8"; Value="-1.26%"}
    @{Ref="D9"; Value="0.9031"}
    @{Ref="E9"; Value="-2.57%"}
    @{Ref="D10"; Value="0.1652"}
    @{Ref="E10"; Value="-6.26%"}
    @{Ref="E11"; Value="-5.38%"}
    @{Ref="D12"; Value="0.08007"}
    @{Ref="E12"; Value="-9.23%"}
    @{Ref="D13"; Value="0.03045"}
    @{Ref="E13"; Value="-3.63%"}
    @{Ref="E14"; Value="-0.40%"}
    @{Ref="D15"; Value="0.001495"}
    @{Ref="E15"; Value="-1.19%"}
    @{Ref="D16"; Value="0.005660"}
    @{Ref="E16"; Value="-5.73%"}
    @{Ref="D17"; Value="3.461"}
    @{Ref="E17"; Value="-0.18%"}
    @{Ref="D18"; Value="2.105"}
    @{Ref="D19"; Value="0.3293"}
    @{Ref="E19"; Value="0.17%"}
    @{Ref="D20"; Value="0.1304"}
    @{Ref="E20"; Value="1.03%"}
    @{Ref="D21"; Value="4.399"}
    @{Ref="E21"; Value="4.61%"}
    @{Ref="D22"; Value="0.2003"}
    @{Ref="E22"; Value="11.81%"}
    @{Ref="D23"; Value="0.04488"}
    @{Ref="E23"; Value="-2.59%"}
    @{Ref="D24"; Value="0.001220"}
    @{Ref="E24"; Value="-1.30%"}
    @{Ref="D25"; Value="0.004024"}
    @{Ref="E25"; Value="-10.44%"}
    @{Ref="D26"; Value="0.0001256"}
    @{Ref="E26"; Value="0.54%"}
    @{Ref="D39"; Value="0.01658"}
    @{Ref="E39"; Value="-4.61%"}
    @{Ref="D40"; Value="0.04347"}
    @{Ref="E40"; Value="-9.58%"}
    @{Ref="D41"; Value="0.007460"}
    @{Ref="E41"; Value="1.55%"}
    @{Ref="D42"; Value="0.1316"}
    @{Ref="E42"; Value="-3.79%"}
    @{Ref="D43"; Value="0.002048"}
    @{Ref="E43"; Value="-12.47%"}
    @{Ref="D44"; Value="0.01022"}
    @{Ref="E44"; Value="-7.11%"}
    @{Ref="D45"; Value="0.00005719"}
    @{Ref="E45"; Value="-5.44%"}
    @{Ref="D46"; Value="0.00000000754"}
    @{Ref="E46"; Value="0.52%"}
    @{Ref="D47"; Value="2.175"}
    @{Ref="E47"; Value="165.12%"}
    @{Ref="D48"; Value="0.003015"}
    @{Ref="E48"; Value="-11.04%"}
    @{Ref="D49"; Value="0.00002110"}
    @{Ref="E49"; Value="0.52%"}
    @{Ref="D50"; Value="0.0002010"}
    @{Ref="E50"; Value="0.52%"}
)

foreach ($update in $updates) {
    $rng = $ws.Range($update.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $update.Value
    $rng.Style = "Normal"
}
